$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "中际旭创"
$ws.Range("B2").Value = "中际旭创"
$ws.Range("C2").Value = "航天发展"
$ws.Range("A3").Value = "实达集团"
$ws.Range("B3").Value = "航天发展"
$ws.Range("C3").Value = "实达集团"
$ws.Range("A4").Value = "特发信息"
$ws.Range("B4").Value = "实达集团"
$ws.Range("C4").Value = "东百集团"
$ws.Range("A5").Value = "航天发展"
$ws.Range("B5").Value = "永鼎股份"
$ws.Range("C5").Value = "欢瑞世纪"
$ws.Range("A6").Value = "雷科防务"
$ws.Range("B6").Value = "雷科防务"
$ws.Range("C6").Value = "特发信息"
$ws.Range("A7").Value = "东百集团"
$ws.Range("B7").Value = "特发信息"
$ws.Range("C7").Value = "中际旭创"
$ws.Range("A8").Value = "新 华 都"
$ws.Range("B8").Value = "海南海药"
$ws.Range("C8").Value = "雷科防务"
$ws.Range("A9").Value = "永鼎股份"
$ws.Range("B9").Value = "蓝色光标"
$ws.Range("C9").Value = "遥望科技"
$ws.Range("A10").Value = "赛微电子"
$ws.Range("B10").Value = "榕基软件"
$ws.Range("C10").Value = "平潭发展"
$ws.Range("A11").Value = "欢瑞世纪"
$ws.Range("B11").Value = "工业富联"
$ws.Range("C11").Value = "永鼎股份"
$ws.Range("A12").Value = "榕基软件"
$ws.Range("B12").Value = "赛微电子"
$ws.Range("C12").Value = "榕基软件"
$ws.Range("A13").Value = "光库科技"
$ws.Range("B13").Value = "东百集团"
$ws.Range("C13").Value = "三江购物"
$ws.Range("A14").Value = "平潭发展"
$ws.Range("B14").Value = "平潭发展"
$ws.Range("C14").Value = "国晟科技"
$ws.Range("A15").Value = "新易盛"
$ws.Range("B15").Value = "欢瑞世纪"
$ws.Range("C15").Value = "石基信息"
$ws.Range("A16").Value = "海南海药"
$ws.Range("B16").Value = "新易盛"
$ws.Range("C16").Value = "赛微电子"
$ws.Range("A17").Value = "蓝色光标"
$ws.Range("B17").Value = "沪电股份"
$ws.Range("C17").Value = "新华都"
$ws.Range("A18").Value = "工业富联"
$ws.Range("B18").Value = "光库科技"
$ws.Range("C18").Value = "合富中国"
$ws.Range("A19").Value = "航天动力"
$ws.Range("B19").Value = "大洋电机"
$ws.Range("C19").Value = "国光连锁"
$ws.Range("A20").Value = "国晟科技"
$ws.Range("B20").Value = "新 华 都"
$ws.Range("C20").Value = "万科A"
$ws.Range("A21").Value = "三江购物"
$ws.Range("B21").Value = "达华智能"
$ws.Range("C21").Value = "蓝色光标"